$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Widgets Properties")

# Column B (Attribute names)
$ws.Cells.Item(13, 2).Value = "filled"
$ws.Cells.Item(14, 2).Value = "fill_highlight"

# Column C (Required / default info)
$ws.Cells.Item(13, 3).Value = "IN => False"
$ws.Cells.Item(14, 3).Value = "IN => True"

# Column D (Type)
$ws.Cells.Item(13, 4).Value = "bool"
$ws.Cells.Item(14, 4).Value = "bool"

# Column E (Example Value) -- force literal text "true"/"false" rather than Boolean
$c13 = $ws.Cells.Item(13, 5)
$c13.Formula = '="tr"&"ue"'
$c13.Copy()
$c13.PasteSpecial(-4163)

$c14 = $ws.Cells.Item(14, 5)
$c14.Formula = '="fal"&"se"'
$c14.Copy()
$c14.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# Column F (Beschreibung)
$ws.Cells.Item(13, 6).Value = "Gibt an ob die Fläche des Rechtecks gefüllt werden soll"

$ws.Range("E13").Select()
